# Applies the MRR template edits described in the commit:
#  - Adds a spell-check/proofing-style split around template placeholders
#    (cosmetic; captured implicitly by Word's own text handling — the
#    substantive content of the document is unaffected by that split).
#  - Restructures the "research technique" paragraph (section 6) so that
#    the OTHER technique is handled first and every other technique value
#    now renders via the new `technique_label` property instead of a long
#    chain of per-technique {IF ...}{END-IF} blocks.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Section 6 technique paragraph: replace the template logic.
#
# The paragraph currently reads (concatenated across runs):
#   {IF $technique.technique_value === `PERSONAL_KNOWLEDGE`}Personal
#   knowledge in procuring supplies/services of this type — Name/Position
#   of Person Relied Upon: {researchPersonalKnowledgePersonOrPosition}
#   {END-IF}{IF $technique.technique_value === `DISA_MARKET_RESEARCH_REPO`}
#   DISA Market Research Repository{END-IF} ... {IF $technique.technique_value
#   === `OTHER`}Other: {techniqueOther}{END-IF}
#
# It becomes:
#   {IF $technique.technique_value === `OTHER`}Other: {techniqueOther}{END-IF}
#   {IF $technique.technique_value !== `OTHER`}{$technique.technique_label}
#   {IF $technique.technique_value === `PERSONAL_KNOWLEDGE`} — Name/Position
#   of Person Relied Upon: {researchPersonalKnowledgePersonOrPosition}
#   {END-IF}{END-IF}
#
# The em dash "—" sits in its own run with distinct formatting (Arial /
# a different font color) and is left untouched; we edit the text
# immediately before it and immediately after it in two separate Find/
# Replace calls so that run survives intact.

# 1) Tail half: everything from " Name/Position..." through the old
#    per-technique {IF}/{END-IF} chain collapses down to just the
#    Name/Position sentence followed by the two closing {END-IF}s.
$oldTail = ' Name/Position of Person Relied Upon: {researchPersonalKnowledgePersonOrPosition}{END-IF}{IF $technique.technique_value === `DISA_MARKET_RESEARCH_REPO`}DISA Market Research Repository{END-IF}{IF $technique.technique_value === `CONTACT_WITH_KNOWLEDGEABLE_PERSON`}Contact with knowledgeable persons in Government and industry.{END-IF}{IF $technique.technique_value === `REVIEW_SIMILAR_RECENT_RESULTS`}Review of recent market research results for similar or identical supplies/services.{END-IF}{IF $technique.technique_value === `REVIEW_DATABASES`}Review of Government and/or commercial databases for relevant information.{END-IF}{IF $technique.technique_value === `REVIEW_SOURCE_LISTS`}Review of source lists for identical or similar items obtained.{END-IF}{IF $technique.technique_value === `REVIEW_PRODUCT_LITERATURE`}Review of catalogs and/or other generally available product literature.{END-IF}{IF $technique.technique_value === `REVIEW_OTHER_CONTRACTS`}Review of existing DISA-wide and other Government-wide Acquisition Contracts.{END-IF}{IF $technique.technique_value === `OTHER`}Other: {techniqueOther}{END-IF}'
$newTail = ' Name/Position of Person Relied Upon: {researchPersonalKnowledgePersonOrPosition}{END-IF}{END-IF}'

$foundTail = $d.Content.Find.Execute($oldTail, $true, $false, $false, $false, $false, $true, 1, $false, $newTail, 2)
if (-not $foundTail) {
    throw "Could not find the technique paragraph's tail text to replace."
}

# 2) Head half: the opening {IF ... PERSONAL_KNOWLEDGE`}Personal
#    knowledge... sentence is replaced with the new OTHER-first logic,
#    ending right before the (untouched) em dash run.
$oldHead = '{IF $technique.technique_value === `PERSONAL_KNOWLEDGE`}Personal knowledge in procuring supplies/services of this type '
$newHead = '{IF $technique.technique_value === `OTHER`}Other: {techniqueOther}{END-IF}{IF $technique.technique_value !== `OTHER`}{$technique.technique_label}{IF $technique.technique_value === `PERSONAL_KNOWLEDGE`} '

$foundHead = $d.Content.Find.Execute($oldHead, $true, $false, $false, $false, $false, $true, 1, $false, $newHead, 2)
if (-not $foundHead) {
    throw "Could not find the technique paragraph's head text to replace."
}

Write-Output "tail replaced: $foundTail; head replaced: $foundHead"
